$wb = $excel.ActiveWorkbook

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1318.25
$ws.Range("I28").Value = 1356.2727
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 1356.2727
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -871.2727
$ws.Range("N28").Value = -1870

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2979
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2979
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2979
$ws.Range("N51").Value = -3947
$ws.Range("M51").ClearContents()

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5200.727
$ws.Range("I70").Value = 3954
$ws.Range("J70").Value = 5913.143
$ws.Range("K70").Value = 11862
$ws.Range("L70").Value = 17739.429
$ws.Range("M70").Value = -11592
$ws.Range("N70").Value = -18279.429

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5200.727
$ws.Range("I73").Value = 3954
$ws.Range("J73").Value = 5913.143
$ws.Range("K73").Value = 11862
$ws.Range("L73").Value = 17739.429
$ws.Range("M73").Value = -10926
$ws.Range("N73").Value = -19611.429

# Sheet ALC, row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 67619.12
$ws.Range("I74").Value = 114218.71
$ws.Range("J74").Value = 13252.917
$ws.Range("K74").Value = 114218.71
$ws.Range("L74").Value = 13252.917
$ws.Range("M74").Value = -113282.71
$ws.Range("N74").Value = -15124.917

# Sheet ALC, row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 67619.12
$ws.Range("I77").Value = 114218.71
$ws.Range("J77").Value = 13252.917
$ws.Range("K77").Value = 571093.55
$ws.Range("L77").Value = 66264.58499999999
$ws.Range("M77").Value = -566413.55
$ws.Range("N77").Value = -75624.58499999999

# Sheet ALC, row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1045
$ws.Range("I103").Value = 663.25
$ws.Range("J103").Value = 1299.5
$ws.Range("K103").Value = 1989.75
$ws.Range("L103").Value = 3898.5
$ws.Range("M103").Value = -1403.75
$ws.Range("N103").Value = -5070.5

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3104.8333
$ws.Range("I113").Value = 2928.6667
$ws.Range("J113").Value = 3633.3333
$ws.Range("K113").Value = 2928.6667
$ws.Range("L113").Value = 3633.3333
$ws.Range("M113").Value = 325.3332999999998
$ws.Range("N113").Value = -10141.3333

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10666
$ws.Range("I116").Value = 14998
$ws.Range("J116").Value = 8500
$ws.Range("K116").Value = 14998
$ws.Range("L116").Value = 8500
$ws.Range("M116").Value = -11556
$ws.Range("N116").Value = -15384

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2117.75
$ws.Range("I137").Value = 2134.5715
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 6403.7145
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3853.7145
$ws.Range("N137").Value = -11100

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4088.1228
$ws.Range("I138").Value = 3511.6924
$ws.Range("J138").Value = 4258.4316
$ws.Range("K138").Value = 10535.0772
$ws.Range("L138").Value = 12775.2948
$ws.Range("M138").Value = -5395.0772
$ws.Range("N138").Value = -23055.2948

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1677.4
$ws.Range("I141").Value = 1595.8334
$ws.Range("J141").Value = 1799.75
$ws.Range("K141").Value = 4787.5002
$ws.Range("L141").Value = 5399.25
$ws.Range("M141").Value = 392.4997999999996
$ws.Range("N141").Value = -15759.25

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8598.682000000001
$ws.Range("I32").Value = 6865.3335
$ws.Range("J32").Value = 44999
$ws.Range("K32").Value = 6865.3335
$ws.Range("L32").Value = 44999
$ws.Range("M32").Value = -6578.3335
$ws.Range("N32").Value = -45573

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3557.5417
$ws.Range("I61").Value = 1392.6666
$ws.Range("J61").Value = 7165.6665
$ws.Range("K61").Value = 1392.6666
$ws.Range("L61").Value = 7165.6665
$ws.Range("M61").Value = -1180.6666
$ws.Range("N61").Value = -7589.6665

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2819
$ws.Range("I74").Value = 2080.9167
$ws.Range("J74").Value = 4295.1665
$ws.Range("K74").Value = 2080.9167
$ws.Range("L74").Value = 4295.1665
$ws.Range("M74").Value = -1206.9167
$ws.Range("N74").Value = -6043.1665

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2819
$ws.Range("I77").Value = 2080.9167
$ws.Range("J77").Value = 4295.1665
$ws.Range("K77").Value = 10404.5835
$ws.Range("L77").Value = 21475.8325
$ws.Range("M77").Value = -6036.583500000001
$ws.Range("N77").Value = -30211.8325

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1112.4286
$ws.Range("I102").Value = 1112.4286
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1112.4286
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 509.5714

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1391.8334
$ws.Range("I132").Value = 1315.875
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 3947.625
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -1417.625
$ws.Range("N132").Value = -11058.5

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3557.5417
$ws.Range("I136").Value = 1392.6666
$ws.Range("J136").Value = 7165.6665
$ws.Range("K136").Value = 4177.9998
$ws.Range("L136").Value = 21496.9995
$ws.Range("M136").Value = -1627.9998
$ws.Range("N136").Value = -26596.9995

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10413.167
$ws.Range("I94").Value = 10413.167
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 10413.167
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -9962.166999999999

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2409.5
$ws.Range("I105").Value = 2425.7896
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 2425.7896
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = -678.7896000000001
$ws.Range("N105").Value = -5594

# Sheet BSM, row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 84499.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 84499.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 84499.5
$ws.Range("N132").Value = -94619.5

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1687.7368
$ws.Range("I16").Value = 1330.3334
$ws.Range("J16").Value = 3028
$ws.Range("K16").Value = 1330.3334
$ws.Range("L16").Value = 3028
$ws.Range("M16").Value = -1043.3334
$ws.Range("N16").Value = -3602

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 125002
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 125002
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 125002
$ws.Range("N22").Value = -125702
$ws.Range("M22").ClearContents()

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2549
$ws.Range("I99").Value = 1998.25
$ws.Range("J99").Value = 3099.75
$ws.Range("K99").Value = 1998.25
$ws.Range("L99").Value = 3099.75
$ws.Range("M99").Value = -500.25
$ws.Range("N99").Value = -6095.75

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1687.7368
$ws.Range("I113").Value = 1330.3334
$ws.Range("J113").Value = 3028
$ws.Range("K113").Value = 1330.3334
$ws.Range("L113").Value = 3028
$ws.Range("M113").Value = 839.6666
$ws.Range("N113").Value = -7368

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2158.5557
$ws.Range("I122").Value = 1933
$ws.Range("J122").Value = 2271.3333
$ws.Range("K122").Value = 5799
$ws.Range("L122").Value = 6813.999899999999
$ws.Range("M122").Value = -3349
$ws.Range("N122").Value = -11713.9999

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2549
$ws.Range("I126").Value = 1998.25
$ws.Range("J126").Value = 3099.75
$ws.Range("K126").Value = 5994.75
$ws.Range("L126").Value = 9299.25
$ws.Range("M126").Value = -3524.75
$ws.Range("N126").Value = -14239.25

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2355.818
$ws.Range("I132").Value = 1273.7142
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 3821.1426
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -1291.1426
$ws.Range("N132").Value = -17808.5

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3821.5454
$ws.Range("I134").Value = 3703.7
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 11111.1
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -8576.099999999999
$ws.Range("N134").Value = -20070

# Sheet CUL, row 57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10828.333
$ws.Range("I57").Value = 3491
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 10473
$ws.Range("L57").Value = 60000
$ws.Range("M57").Value = -9914
$ws.Range("N57").Value = -61118

# Sheet CUL, row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 17252.166
$ws.Range("I76").Value = 11256.5
$ws.Range("J76").Value = 20250
$ws.Range("K76").Value = 33769.5
$ws.Range("L76").Value = 60750
$ws.Range("M76").Value = -33386.5
$ws.Range("N76").Value = -61516

# Sheet CUL, row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 17252.166
$ws.Range("I79").Value = 11256.5
$ws.Range("J79").Value = 20250
$ws.Range("K79").Value = 33769.5
$ws.Range("L79").Value = 60750
$ws.Range("M79").Value = -32443.5
$ws.Range("N79").Value = -63402

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 7971.467
$ws.Range("I140").Value = 3286.111
$ws.Range("J140").Value = 14999.5
$ws.Range("K140").Value = 9858.332999999999
$ws.Range("L140").Value = 44998.5
$ws.Range("M140").Value = -4678.332999999999
$ws.Range("N140").Value = -55358.5

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1525.4117
$ws.Range("I102").Value = 1308.25
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1308.25
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 313.75
$ws.Range("N102").Value = -8244

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3963.7144
$ws.Range("I107").Value = 1081.3334
$ws.Range("J107").Value = 6125.5
$ws.Range("K107").Value = 1081.3334
$ws.Range("L107").Value = 6125.5
$ws.Range("M107").Value = 838.6666
$ws.Range("N107").Value = -9965.5

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9630.933000000001
$ws.Range("I46").Value = 85112.164
$ws.Range("J46").Value = 1085.8868
$ws.Range("K46").Value = 85112.164
$ws.Range("L46").Value = 1085.8868
$ws.Range("M46").Value = -84924.164
$ws.Range("N46").Value = -1461.8868

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3338.0527
$ws.Range("I61").Value = 2690.2778
$ws.Range("J61").Value = 14998
$ws.Range("K61").Value = 2690.2778
$ws.Range("L61").Value = 14998
$ws.Range("M61").Value = -2488.2778
$ws.Range("N61").Value = -15402

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4624.875
$ws.Range("I100").Value = 3666.3333
$ws.Range("J100").Value = 5200
$ws.Range("K100").Value = 3666.3333
$ws.Range("L100").Value = 5200
$ws.Range("M100").Value = -3125.3333
$ws.Range("N100").Value = -6282

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3338.0527
$ws.Range("I113").Value = 2690.2778
$ws.Range("J113").Value = 14998
$ws.Range("K113").Value = 2690.2778
$ws.Range("L113").Value = 14998
$ws.Range("M113").Value = -520.2777999999998
$ws.Range("N113").Value = -19338

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 855.1579
$ws.Range("I132").Value = 897.8
$ws.Range("J132").Value = 695.25
$ws.Range("K132").Value = 2693.4
$ws.Range("L132").Value = 2085.75
$ws.Range("M132").Value = -163.3999999999996
$ws.Range("N132").Value = -7145.75

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -12600
$ws.Range("M136").ClearContents()

# Sheet WVR, row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1999
$ws.Range("I14").Value = 1999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1831
$ws.Range("N14").ClearContents()

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 10000506
$ws.Range("I100").Value = 20000372
$ws.Range("J100").Value = 640
$ws.Range("K100").Value = 40000744
$ws.Range("L100").Value = 1280
$ws.Range("M100").Value = -40000203
$ws.Range("N100").Value = -2362

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3422.5715
$ws.Range("I132").Value = 3493.75
$ws.Range("J132").Value = 3327.6667
$ws.Range("K132").Value = 10481.25
$ws.Range("L132").Value = 9983.000100000001
$ws.Range("M132").Value = -7951.25
$ws.Range("N132").Value = -15043.0001

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2564.0908
$ws.Range("I136").Value = 3129.1428
$ws.Range("J136").Value = 1575.25
$ws.Range("K136").Value = 9387.428400000001
$ws.Range("L136").Value = 4725.75
$ws.Range("M136").Value = -6837.428400000001
$ws.Range("N136").Value = -9825.75
